$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Fal files" / "Fal images" columns (headers)
$ws.Range("S1").Value = "Fal files"
$ws.Range("T1").Value = "Fal images"

# Row 2 (Chevrolet Pick up) - Fal images
$ws.Range("T2").Value = "bicycle_cannondale_black.jpg,bicycle_cannondale_red.jpg,bicycle_cannondale_white.jpg"

# Row 3 (BMW Roadster) - Fal files
$ws.Range("S3").Value = "pickup_chevrolet_black.jpg,pickup_chevrolet_red.jpg,pickup_chevrolet_white.jpg"

# Row 4 (Cannondale) - Fal files / Fal images
$ws.Range("S4").Value = "bicycle_cannondale_black.jpg"
$ws.Range("T4").Value = "roadster_bmw_red.jpg,roadster_bmw_black.jpg"

# Update selection to match the authored state (active cell moves to T4
# after entering the last value, and the view scrolls right to column E)
$ws.Range("T4").Select()
